$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169794201850891
$ws.Range("B1").Value = 2.161981344223022
$ws.Range("C1").Value = 3.324399471282959
$ws.Range("D1").Value = 3.660769701004028
$ws.Range("E1").Value = 1.172302007675171
